$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates (coinranking.com crypto table refresh) ---
# D-column (Price) cells are plain text in the source sheet (e.g. "59.527.78" uses
# "." as a thousands separator), so we force text format before assigning the new
# value to avoid Excel re-interpreting it as a locale-parsed number, then restore
# the cell style to "Normal" so no stray formatting is left behind.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "59.527.78"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.20%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.995.40"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.44%  "

$ws.Range("E4").Value = "  +0.08%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "563.76"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.86%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "138.76"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.82%  "

$ws.Range("E7").Value = "  -0.10%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.520"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.57%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.981.35"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.23%  "

$ws.Range("E10").Value = "  +5.24%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +11.55%  "

$ws.Range("E12").Value = "  +2.13%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000229"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.17%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "33.79"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.59%  "

$ws.Range("E15").Value = "  -0.21%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.491.34"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.50%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "7.17"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.64%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.991.19"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.42%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "59.511.96"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.21%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "434.98"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.71%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.56"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.717"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.97%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "13.49"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "7.06"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "80.19"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").Value = "  +11.06%  "

$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("E29").Value = "  +3.26%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.76"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.92%  "

$ws.Range("E31").Value = "  +5.09%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "25.76"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.105"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +9.46%  "

$ws.Range("E34").Value = "  +13.49%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.08"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.63%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "48.90"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("E39").Value = "  -1.81%  "

$ws.Range("E40").Value = "  +7.84%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "402.84"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +8.61%  "

$ws.Range("E42").Value = "  +3.00%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.761.17"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.23%  "

$ws.Range("E44").Value = "  -1.51%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.250"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.92%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "123.11"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "34.49"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +19.50%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "23.55"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.86%  "

# --- Row 35/36 swap: Filecoin <-> Mantle, each with refreshed price/volume figures ---
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.986"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.22%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.87"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.84%  "

# --- Row 49/50 swap: Fetch.AI <-> Stellar, each with refreshed price/volume figures ---
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.110"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.05%  "

$ws.Range("B50").Value = "Fetch.AI"
$ws.Range("C50").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.98%  "

